$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 218.74359  # H17: 173.06522 -> 218.74359
$ws.Cells.Item(17, 10).Value = 216.13889  # J17: 167.69768 -> 216.13889
$ws.Cells.Item(17, 12).Value = 648.4166700000001  # L17: 503.09304 -> 648.4166700000001
$ws.Cells.Item(17, 14).Value = -984.4166700000001  # N17: -839.09304 -> -984.4166700000001

$ws.Cells.Item(29, 8).Value = 633.3333  # H29: 800 -> 633.3333
$ws.Cells.Item(29, 10).Value = 1500  # J29: 2000 -> 1500
$ws.Cells.Item(29, 12).Value = 4500  # L29: 6000 -> 4500
$ws.Cells.Item(29, 14).Value = -5062  # N29: -6562 -> -5062

$ws.Cells.Item(31, 8).Value = 2500  # H31: 2281.5 -> 2500
$ws.Cells.Item(31, 9).Value = 0  # I31: 589 -> 0
$ws.Cells.Item(31, 10).Value = 2500  # J31: 2620 -> 2500
$ws.Cells.Item(31, 11).Value = 0  # K31: 1767 -> 0
$ws.Cells.Item(31, 12).ClearContents()  # L31: 7860 -> (removed)
$ws.Cells.Item(31, 13).Value = 7500  # M31: -1537 -> 7500
$ws.Cells.Item(31, 14).Value = -7960  # N31: -8320 -> -7960

$ws.Cells.Item(115, 8).Value = 3694  # H115: 3351 -> 3694
$ws.Cells.Item(115, 9).Value = 5235  # I115: 3585 -> 5235
$ws.Cells.Item(115, 10).Value = 2666.6667  # J115: 3000 -> 2666.6667
$ws.Cells.Item(115, 11).Value = 15705  # K115: 10755 -> 15705
$ws.Cells.Item(115, 12).Value = 8000.000100000001  # L115: 9000 -> 8000.000100000001
$ws.Cells.Item(115, 13).Value = -14138  # M115: -9188 -> -14138
$ws.Cells.Item(115, 14).Value = -11134.0001  # N115: -12134 -> -11134.0001

$ws.Cells.Item(121, 8).Value = 2772.3914  # H121: 2743.25 -> 2772.3914
$ws.Cells.Item(121, 9).Value = 1100  # I121: 1200 -> 1100
$ws.Cells.Item(121, 10).Value = 2848.4092  # J121: 2824.4736 -> 2848.4092
$ws.Cells.Item(121, 11).Value = 3300  # K121: 3600 -> 3300
$ws.Cells.Item(121, 12).Value = 8545.2276  # L121: 8473.4208 -> 8545.2276
$ws.Cells.Item(121, 13).Value = -1553  # M121: -1853 -> -1553
$ws.Cells.Item(121, 14).Value = -12039.2276  # N121: -11967.4208 -> -12039.2276

$ws.Cells.Item(127, 8).Value = 975.5714  # H127: 1673.4445 -> 975.5714
$ws.Cells.Item(127, 9).Value = 622.7  # I127: 0 -> 622.7
$ws.Cells.Item(127, 10).Value = 1296.3636  # J127: 1673.4445 -> 1296.3636
$ws.Cells.Item(127, 11).Value = 1868.1  # K127: 0 -> 1868.1
$ws.Cells.Item(127, 12).Value = 3889.0908  # L127: 5020.333500000001 -> 3889.0908
$ws.Cells.Item(127, 13).Value = 3091.9  # M127: None -> 3091.9
$ws.Cells.Item(127, 14).Value = -13809.0908  # N127: -14940.3335 -> -13809.0908

$ws.Cells.Item(129, 8).Value = 1207.7273  # H129: 1177.4474 -> 1207.7273
$ws.Cells.Item(129, 10).Value = 1436.4783  # J129: 1354.5358 -> 1436.4783
$ws.Cells.Item(129, 12).Value = 4309.4349  # L129: 4063.6074 -> 4309.4349
$ws.Cells.Item(129, 14).Value = -14309.4349  # N129: -14063.6074 -> -14309.4349

$ws.Cells.Item(132, 8).Value = 4044.303  # H132: 4804.2964 -> 4044.303
$ws.Cells.Item(132, 9).Value = 3985.8333  # I132: 4708.68 -> 3985.8333
$ws.Cells.Item(132, 10).Value = 4629  # J132: 5999.5 -> 4629
$ws.Cells.Item(132, 11).Value = 11957.4999  # K132: 14126.04 -> 11957.4999
$ws.Cells.Item(132, 12).Value = 13887  # L132: 17998.5 -> 13887
$ws.Cells.Item(132, 13).Value = -9427.499899999999  # M132: -11596.04 -> -9427.499899999999
$ws.Cells.Item(132, 14).Value = -18947  # N132: -23058.5 -> -18947

$ws.Cells.Item(133, 8).Value = 80280  # H133: 60780 -> 80280
$ws.Cells.Item(133, 10).Value = 80280  # J133: 60780 -> 80280
$ws.Cells.Item(133, 12).Value = 80280  # L133: 60780 -> 80280
$ws.Cells.Item(133, 14).Value = -90400  # N133: -70900 -> -90400

$ws.Cells.Item(138, 8).Value = 2738.426  # H138: 2697.5964 -> 2738.426
$ws.Cells.Item(138, 10).Value = 2440.5454  # J138: 2410.0425 -> 2440.5454
$ws.Cells.Item(138, 12).Value = 7321.6362  # L138: 7230.127500000001 -> 7321.6362
$ws.Cells.Item(138, 14).Value = -17601.6362  # N138: -17510.1275 -> -17601.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 451  # H5: 471 -> 451
$ws.Cells.Item(5, 9).Value = 0  # I5: 470 -> 0
$ws.Cells.Item(5, 10).Value = 451  # J5: 472 -> 451
$ws.Cells.Item(5, 11).Value = 0  # K5: 470 -> 0
$ws.Cells.Item(5, 12).ClearContents()  # L5: 472 -> (removed)
$ws.Cells.Item(5, 13).Value = 451  # M5: -358 -> 451
$ws.Cells.Item(5, 14).Value = -675  # N5: -696 -> -675

$ws.Cells.Item(7, 8).Value = 41000  # H7: 40980 -> 41000
$ws.Cells.Item(7, 10).Value = 41000  # J7: 40980 -> 41000
$ws.Cells.Item(7, 12).Value = 41000  # L7: 40980 -> 41000
$ws.Cells.Item(7, 14).Value = -41228  # N7: -41208 -> -41228

$ws.Cells.Item(32, 8).Value = 644244.2  # H32: 840832.75 -> 644244.2
$ws.Cells.Item(32, 9).Value = 820700.4399999999  # I32: 1181428.8 -> 820700.4399999999
$ws.Cells.Item(32, 11).Value = 820700.4399999999  # K32: 1181428.8 -> 820700.4399999999
$ws.Cells.Item(32, 13).Value = -820413.4399999999  # M32: -1181141.8 -> -820413.4399999999

$ws.Cells.Item(88, 8).Value = 2486.5908  # H88: 2475.6843 -> 2486.5908
$ws.Cells.Item(88, 9).Value = 2517.7144  # I88: 2474.3333 -> 2517.7144
$ws.Cells.Item(88, 10).Value = 2432.125  # J88: 2500 -> 2432.125
$ws.Cells.Item(88, 11).Value = 2517.7144  # K88: 2474.3333 -> 2517.7144
$ws.Cells.Item(88, 12).Value = 2432.125  # L88: 2500 -> 2432.125
$ws.Cells.Item(88, 13).Value = -2111.7144  # M88: -2068.3333 -> -2111.7144
$ws.Cells.Item(88, 14).Value = -3244.125  # N88: -3312 -> -3244.125

$ws.Cells.Item(91, 8).Value = 2486.5908  # H91: 2475.6843 -> 2486.5908
$ws.Cells.Item(91, 9).Value = 2517.7144  # I91: 2474.3333 -> 2517.7144
$ws.Cells.Item(91, 10).Value = 2432.125  # J91: 2500 -> 2432.125
$ws.Cells.Item(91, 11).Value = 2517.7144  # K91: 2474.3333 -> 2517.7144
$ws.Cells.Item(91, 12).Value = 2432.125  # L91: 2500 -> 2432.125
$ws.Cells.Item(91, 13).Value = -1113.7144  # M91: -1070.3333 -> -1113.7144
$ws.Cells.Item(91, 14).Value = -5240.125  # N91: -5308 -> -5240.125

$ws.Cells.Item(101, 8).Value = 79602  # H101: 75373.14 -> 79602
$ws.Cells.Item(101, 10).Value = 79602  # J101: 75373.14 -> 79602
$ws.Cells.Item(101, 12).Value = 79602  # L101: 75373.14 -> 79602
$ws.Cells.Item(101, 14).Value = -86092  # N101: -81863.14 -> -86092

$ws.Cells.Item(133, 8).Value = 65261  # H133: 50000 -> 65261
$ws.Cells.Item(133, 10).Value = 65261  # J133: 50000 -> 65261
$ws.Cells.Item(133, 12).Value = 65261  # L133: 50000 -> 65261
$ws.Cells.Item(133, 14).Value = -70321  # N133: -55060 -> -70321

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 451  # H4: 471 -> 451
$ws.Cells.Item(4, 9).Value = 0  # I4: 470 -> 0
$ws.Cells.Item(4, 10).Value = 451  # J4: 472 -> 451
$ws.Cells.Item(4, 11).Value = 0  # K4: 470 -> 0
$ws.Cells.Item(4, 12).ClearContents()  # L4: 472 -> (removed)
$ws.Cells.Item(4, 13).Value = 451  # M4: -355 -> 451
$ws.Cells.Item(4, 14).Value = -681  # N4: -702 -> -681

$ws.Cells.Item(86, 8).Value = 90911800  # H86: 66669090 -> 90911800
$ws.Cells.Item(86, 9).Value = 142858980  # I86: 125001760 -> 142858980
$ws.Cells.Item(86, 10).Value = 4251.75  # J86: 3173.4285 -> 4251.75
$ws.Cells.Item(86, 11).Value = 142858980  # K86: 125001760 -> 142858980
$ws.Cells.Item(86, 12).Value = 4251.75  # L86: 3173.4285 -> 4251.75
$ws.Cells.Item(86, 13).Value = -142857857  # M86: -125000637 -> -142857857
$ws.Cells.Item(86, 14).Value = -6497.75  # N86: -5419.4285 -> -6497.75

$ws.Cells.Item(89, 8).Value = 90911800  # H89: 66669090 -> 90911800
$ws.Cells.Item(89, 9).Value = 142858980  # I89: 125001760 -> 142858980
$ws.Cells.Item(89, 10).Value = 4251.75  # J89: 3173.4285 -> 4251.75
$ws.Cells.Item(89, 11).Value = 714294900  # K89: 625008800 -> 714294900
$ws.Cells.Item(89, 12).Value = 21258.75  # L89: 15867.1425 -> 21258.75
$ws.Cells.Item(89, 13).Value = -714289284  # M89: -625003184 -> -714289284
$ws.Cells.Item(89, 14).Value = -32490.75  # N89: -27099.1425 -> -32490.75

$ws.Cells.Item(93, 8).Value = 29000  # H93: 29525 -> 29000
$ws.Cells.Item(93, 10).Value = 29000  # J93: 29525 -> 29000
$ws.Cells.Item(93, 12).Value = 29000  # L93: 29525 -> 29000
$ws.Cells.Item(93, 14).Value = -32744  # N93: -33269 -> -32744

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4919.3784  # H31: 5442.606 -> 4919.3784
$ws.Cells.Item(31, 9).Value = 1394.7894  # I31: 1499.3529 -> 1394.7894
$ws.Cells.Item(31, 10).Value = 8639.777  # J31: 9632.3125 -> 8639.777
$ws.Cells.Item(31, 11).Value = 1394.7894  # K31: 1499.3529 -> 1394.7894
$ws.Cells.Item(31, 12).Value = 8639.777  # L31: 9632.3125 -> 8639.777
$ws.Cells.Item(31, 13).Value = -1099.7894  # M31: -1204.3529 -> -1099.7894
$ws.Cells.Item(31, 14).Value = -9229.777  # N31: -10222.3125 -> -9229.777

$ws.Cells.Item(34, 8).Value = 4919.3784  # H34: 5442.606 -> 4919.3784
$ws.Cells.Item(34, 9).Value = 1394.7894  # I34: 1499.3529 -> 1394.7894
$ws.Cells.Item(34, 10).Value = 8639.777  # J34: 9632.3125 -> 8639.777
$ws.Cells.Item(34, 11).Value = 1394.7894  # K34: 1499.3529 -> 1394.7894
$ws.Cells.Item(34, 12).Value = 8639.777  # L34: 9632.3125 -> 8639.777
$ws.Cells.Item(34, 13).Value = -1192.7894  # M34: -1297.3529 -> -1192.7894
$ws.Cells.Item(34, 14).Value = -9043.777  # N34: -10036.3125 -> -9043.777

$ws.Cells.Item(58, 8).Value = 2428.0715  # H58: 2711.5833 -> 2428.0715
$ws.Cells.Item(58, 9).Value = 2149.3  # I58: 2517.375 -> 2149.3
$ws.Cells.Item(58, 10).Value = 3125  # J58: 3100 -> 3125
$ws.Cells.Item(58, 11).Value = 2149.3  # K58: 2517.375 -> 2149.3
$ws.Cells.Item(58, 12).Value = 3125  # L58: 3100 -> 3125
$ws.Cells.Item(58, 13).Value = -1946.3  # M58: -2314.375 -> -1946.3
$ws.Cells.Item(58, 14).Value = -3531  # N58: -3506 -> -3531

$ws.Cells.Item(96, 8).Value = 90000  # H96: 57874.668 -> 90000
$ws.Cells.Item(96, 10).Value = 90000  # J96: 57874.668 -> 90000
$ws.Cells.Item(96, 12).Value = 90000  # L96: 57874.668 -> 90000
$ws.Cells.Item(96, 14).Value = -95492  # N96: -63366.668 -> -95492

$ws.Cells.Item(125, 8).Value = 0  # H125: 98326 -> 0
$ws.Cells.Item(125, 10).Value = 0  # J125: 98326 -> 0
$ws.Cells.Item(125, 12).ClearContents()  # L125: 98326 -> (removed)
$ws.Cells.Item(125, 14).Value = 0  # N125: -103246 -> 0

$ws.Cells.Item(132, 8).Value = 9806511  # H132: 13891848 -> 9806511
$ws.Cells.Item(132, 9).Value = 1422.5  # I132: 1500 -> 1422.5
$ws.Cells.Item(132, 10).Value = 15154741  # J132: 18521964 -> 15154741
$ws.Cells.Item(132, 11).Value = 4267.5  # K132: 4500 -> 4267.5
$ws.Cells.Item(132, 12).Value = 45464223  # L132: 55565892 -> 45464223
$ws.Cells.Item(132, 13).Value = -1737.5  # M132: -1970 -> -1737.5
$ws.Cells.Item(132, 14).Value = -45469283  # N132: -55570952 -> -45469283

$ws.Cells.Item(134, 8).Value = 1066  # H134: 1126.8334 -> 1066
$ws.Cells.Item(134, 9).Value = 834.1905  # I134: 956.5454999999999 -> 834.1905
$ws.Cells.Item(134, 10).Value = 3500  # J134: 3000 -> 3500
$ws.Cells.Item(134, 11).Value = 2502.5715  # K134: 2869.6365 -> 2502.5715
$ws.Cells.Item(134, 12).Value = 10500  # L134: 9000 -> 10500
$ws.Cells.Item(134, 13).Value = 32.42849999999999  # M134: -334.6364999999996 -> 32.42849999999999
$ws.Cells.Item(134, 14).Value = -15570  # N134: -14070 -> -15570

$ws.Cells.Item(136, 8).Value = 2428.0715  # H136: 2711.5833 -> 2428.0715
$ws.Cells.Item(136, 9).Value = 2149.3  # I136: 2517.375 -> 2149.3
$ws.Cells.Item(136, 10).Value = 3125  # J136: 3100 -> 3125
$ws.Cells.Item(136, 11).Value = 6447.900000000001  # K136: 7552.125 -> 6447.900000000001
$ws.Cells.Item(136, 12).Value = 9375  # L136: 9300 -> 9375
$ws.Cells.Item(136, 13).Value = -3897.900000000001  # M136: -5002.125 -> -3897.900000000001
$ws.Cells.Item(136, 14).Value = -14475  # N136: -14400 -> -14475

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 1402.42  # H68: 1442.72 -> 1402.42
$ws.Cells.Item(68, 9).Value = 731.91174  # I68: 782.4286 -> 731.91174
$ws.Cells.Item(68, 10).Value = 1747.8334  # J68: 1798.2616 -> 1747.8334
$ws.Cells.Item(68, 11).Value = 2195.73522  # K68: 2347.2858 -> 2195.73522
$ws.Cells.Item(68, 12).Value = 5243.5002  # L68: 5394.7848 -> 5243.5002
$ws.Cells.Item(68, 13).Value = -1384.73522  # M68: -1536.2858 -> -1384.73522
$ws.Cells.Item(68, 14).Value = -6865.5002  # N68: -7016.7848 -> -6865.5002

$ws.Cells.Item(71, 8).Value = 1402.42  # H71: 1442.72 -> 1402.42
$ws.Cells.Item(71, 9).Value = 731.91174  # I71: 782.4286 -> 731.91174
$ws.Cells.Item(71, 10).Value = 1747.8334  # J71: 1798.2616 -> 1747.8334
$ws.Cells.Item(71, 11).Value = 6587.20566  # K71: 7041.8574 -> 6587.20566
$ws.Cells.Item(71, 12).Value = 15730.5006  # L71: 16184.3544 -> 15730.5006
$ws.Cells.Item(71, 13).Value = -2531.20566  # M71: -2985.8574 -> -2531.20566
$ws.Cells.Item(71, 14).Value = -23842.5006  # N71: -24296.3544 -> -23842.5006

$ws.Cells.Item(107, 8).Value = 1936.4546  # H107: 2138.151 -> 1936.4546
$ws.Cells.Item(107, 9).Value = 380.83334  # I107: 375.53845 -> 380.83334
$ws.Cells.Item(107, 10).Value = 2370.5813  # J107: 2711 -> 2370.5813
$ws.Cells.Item(107, 11).Value = 1142.50002  # K107: 1126.61535 -> 1142.50002
$ws.Cells.Item(107, 12).Value = 7111.743899999999  # L107: 8133 -> 7111.743899999999
$ws.Cells.Item(107, 13).Value = 777.4999800000001  # M107: 793.38465 -> 777.4999800000001
$ws.Cells.Item(107, 14).Value = -10951.7439  # N107: -11973 -> -10951.7439

$ws.Cells.Item(121, 8).Value = 1011.5893  # H121: 974.9403 -> 1011.5893
$ws.Cells.Item(121, 9).Value = 522.63635  # I121: 513.5454999999999 -> 522.63635
$ws.Cells.Item(121, 10).Value = 1131.1111  # J121: 1065.5714 -> 1131.1111
$ws.Cells.Item(121, 11).Value = 1567.90905  # K121: 1540.6365 -> 1567.90905
$ws.Cells.Item(121, 12).Value = 3393.3333  # L121: 3196.7142 -> 3393.3333
$ws.Cells.Item(121, 13).Value = -257.90905  # M121: -230.6364999999998 -> -257.90905
$ws.Cells.Item(121, 14).Value = -6013.3333  # N121: -5816.7142 -> -6013.3333

$ws.Cells.Item(136, 8).Value = 3495.0833  # H136: 3420.75 -> 3495.0833
$ws.Cells.Item(136, 10).Value = 3632  # J136: 3483.3333 -> 3632
$ws.Cells.Item(136, 12).Value = 10896  # L136: 10449.9999 -> 10896
$ws.Cells.Item(136, 14).Value = -21096  # N136: -20649.9999 -> -21096

$ws.Cells.Item(140, 8).Value = 2368.16  # H140: 2372.8948 -> 2368.16
$ws.Cells.Item(140, 9).Value = 1506.5385  # I140: 1332.2727 -> 1506.5385
$ws.Cells.Item(140, 10).Value = 3301.5833  # J140: 3803.75 -> 3301.5833
$ws.Cells.Item(140, 11).Value = 4519.6155  # K140: 3996.8181 -> 4519.6155
$ws.Cells.Item(140, 12).Value = 9904.749899999999  # L140: 11411.25 -> 9904.749899999999
$ws.Cells.Item(140, 13).Value = 660.3845000000001  # M140: 1183.1819 -> 660.3845000000001
$ws.Cells.Item(140, 14).Value = -20264.7499  # N140: -21771.25 -> -20264.7499

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(39, 8).Value = 22500  # H39: 23000 -> 22500
$ws.Cells.Item(39, 10).Value = 22500  # J39: 23000 -> 22500
$ws.Cells.Item(39, 12).Value = 22500  # L39: 23000 -> 22500
$ws.Cells.Item(39, 14).Value = -23564  # N39: -24064 -> -23564

$ws.Cells.Item(48, 8).Value = 10000  # H48: 9950 -> 10000
$ws.Cells.Item(48, 10).Value = 10000  # J48: 9950 -> 10000
$ws.Cells.Item(48, 12).Value = 10000  # L48: 9950 -> 10000
$ws.Cells.Item(48, 14).Value = -10970  # N48: -10920 -> -10970

$ws.Cells.Item(130, 8).Value = 56390  # H130: 57780 -> 56390
$ws.Cells.Item(130, 10).Value = 56390  # J130: 57780 -> 56390
$ws.Cells.Item(130, 12).Value = 56390  # L130: 57780 -> 56390
$ws.Cells.Item(130, 14).Value = -66430  # N130: -67820 -> -66430

$ws.Cells.Item(132, 8).Value = 2402.2173  # H132: 3270.8235 -> 2402.2173
$ws.Cells.Item(132, 9).Value = 2007.5555  # I132: 3034 -> 2007.5555
$ws.Cells.Item(132, 10).Value = 2963.0527  # J132: 3343.6924 -> 2963.0527
$ws.Cells.Item(132, 11).Value = 6022.666499999999  # K132: 9102 -> 6022.666499999999
$ws.Cells.Item(132, 12).Value = 8889.158100000001  # L132: 10031.0772 -> 8889.158100000001
$ws.Cells.Item(132, 13).Value = -3492.666499999999  # M132: -6572 -> -3492.666499999999
$ws.Cells.Item(132, 14).Value = -13949.1581  # N132: -15091.0772 -> -13949.1581

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 0  # H68: 1500 -> 0
$ws.Cells.Item(68, 10).Value = 0  # J68: 1500 -> 0
$ws.Cells.Item(68, 12).ClearContents()  # L68: 1500 -> (removed)
$ws.Cells.Item(68, 14).Value = 0  # N68: -2998 -> 0

$ws.Cells.Item(71, 8).Value = 0  # H71: 1500 -> 0
$ws.Cells.Item(71, 10).Value = 0  # J71: 1500 -> 0
$ws.Cells.Item(71, 12).ClearContents()  # L71: 7500 -> (removed)
$ws.Cells.Item(71, 14).Value = 0  # N71: -14988 -> 0

$ws.Cells.Item(99, 8).Value = 0  # H99: 26500 -> 0
$ws.Cells.Item(99, 10).Value = 0  # J99: 26500 -> 0
$ws.Cells.Item(99, 12).ClearContents()  # L99: 26500 -> (removed)
$ws.Cells.Item(99, 14).Value = 0  # N99: -32490 -> 0

$ws.Cells.Item(130, 8).Value = 0  # H130: 39795 -> 0
$ws.Cells.Item(130, 10).Value = 0  # J130: 39795 -> 0
$ws.Cells.Item(130, 12).ClearContents()  # L130: 39795 -> (removed)
$ws.Cells.Item(130, 14).Value = 0  # N130: -49835 -> 0

$ws.Cells.Item(132, 8).Value = 3200.3157  # H132: 3454.818 -> 3200.3157
$ws.Cells.Item(132, 9).Value = 2869.08  # I132: 3206.2 -> 2869.08
$ws.Cells.Item(132, 11).Value = 8607.24  # K132: 9618.599999999999 -> 8607.24
$ws.Cells.Item(132, 13).Value = -6077.24  # M132: -7088.599999999999 -> -6077.24

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 40407  # H46: 50804.668 -> 40407
$ws.Cells.Item(46, 10).Value = 40407  # J46: 50804.668 -> 40407
$ws.Cells.Item(46, 12).Value = 40407  # L46: 50804.668 -> 40407
$ws.Cells.Item(46, 14).Value = -40869  # N46: -51266.668 -> -40869

$ws.Cells.Item(110, 8).Value = 74500  # H110: 200000 -> 74500
$ws.Cells.Item(110, 10).Value = 74500  # J110: 200000 -> 74500
$ws.Cells.Item(110, 12).Value = 74500  # L110: 200000 -> 74500
$ws.Cells.Item(110, 14).Value = -82680  # N110: -208180 -> -82680

$ws.Cells.Item(123, 8).Value = 25143  # H123: 24714 -> 25143
$ws.Cells.Item(123, 10).Value = 25143  # J123: 24714 -> 25143
$ws.Cells.Item(123, 12).Value = 25143  # L123: 24714 -> 25143
$ws.Cells.Item(123, 14).Value = -34943  # N123: -34514 -> -34943

$ws.Cells.Item(130, 8).Value = 30000  # H130: 64214.5 -> 30000
$ws.Cells.Item(130, 10).Value = 30000  # J130: 64214.5 -> 30000
$ws.Cells.Item(130, 12).Value = 30000  # L130: 64214.5 -> 30000
$ws.Cells.Item(130, 14).Value = -40040  # N130: -74254.5 -> -40040

$ws.Cells.Item(134, 8).Value = 40407  # H134: 50804.668 -> 40407
$ws.Cells.Item(134, 10).Value = 40407  # J134: 50804.668 -> 40407
$ws.Cells.Item(134, 12).Value = 121221  # L134: 152414.004 -> 121221
$ws.Cells.Item(134, 14).Value = -126291  # N134: -157484.004 -> -126291
